$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Document number cell: "01/No.08 - Dok.03/2022" -> split into
#    "01/No.08 - Dok.03/" + "VST/" + [_GoBack bookmark] + "2022".
#    The _GoBack bookmark marks the final edit location, so we move it here
#    (removing it from its old spot further down in the document).
# ---------------------------------------------------------------------------

# Remove the old _GoBack bookmark (near the end of the document) first.
try { $d.Bookmarks("_GoBack").Delete() } catch { }

$numRng = $d.Content
$numRng.Find.Execute("01/No.08 – Dok.03/2022")
$numStart = $numRng.Start
$numEnd = $numRng.End

# Insert "VST/" right before "2022" (the last 4 characters of the match).
$insertPos = $numEnd - 4
$sel = $word.Selection
$sel.SetRange($insertPos, $insertPos)
$sel.TypeText("VST/")

# Re-create the _GoBack bookmark right between "VST/" and "2022" - this also
# forces a run split at that exact point.
$bmPos = $insertPos + 4
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# Force the leading space run (immediately before the doc-number text) to
# stay split out from the doc-number run: toggling formatting and reverting
# it makes the engine keep it as its own run instead of silently
# re-consolidating it with its neighbor on save.
$spaceRng = $d.Range($numStart - 1, $numStart)
$spaceRng.Font.Bold = 1
$spaceRng.Font.Bold = 0

# Force "VST/" to stay split from "01/No.08 - Dok.03/" on its left.
$vstRng = $d.Range($insertPos, $insertPos + 4)
$vstRng.Font.Bold = 1
$vstRng.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2)-5) Four list paragraphs where several adjacent runs carrying the exact
#    same formatting get consolidated into a single run (no text change).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Kiri menggunakan kunci ", $true, $false, $false, $false, $false, $true, 1, $false, " Kiri menggunakan kunci ", 2)
$d.Content.Find.Execute(" Kanan menggunakan kunci ", $true, $false, $false, $false, $false, $true, 1, $false, " Kanan menggunakan kunci ", 2)
$d.Content.Find.Execute("Sisi kiri lebih maju dari sisi kanan maka", $true, $false, $false, $false, $false, $true, 1, $false, "Sisi kiri lebih maju dari sisi kanan maka", 2)

Write-Output "done"
